$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.363.15'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '1.825.64'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'314.98"
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = "'0.4476"
$ws.Range('E7').Value = '  -2.13%  '
$ws.Range('D8').Value = "'0.3787"
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').Value = "'0.07478"
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('D10').Value = "'0.8873"
$ws.Range('E10').Value = '  +3.04%  '
$ws.Range('D11').Value = "'21.05"
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.827.08'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = "'6.755"
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = "'5.462"
$ws.Range('D15').Value = "'93.90"
$ws.Range('D16').Value = "'0.07120"
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = "'1.002"
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').Value = "'0.000008791"
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = "'1.001"
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = "'15.17"
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = '27.380.51'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = "'5.408"
$ws.Range('E22').Value = '  +4.13%  '
$ws.Range('D23').Value = "'10.99"
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = "'1.961"
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = "'151.64"
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = "'2.306"
$ws.Range('E26').Value = '  +3.71%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'18.69"
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = "'5.414"
$ws.Range('E28').Value = '  +2.73%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = "'117.83"
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = "'0.08891"
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'0.7918"
$ws.Range('E31').Value = '  +2.64%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = "'1.208"
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'4.613"
$ws.Range('E33').Value = '  +3.27%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = "'2.925"
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('B35').Value = 'Frax'
$ws.Range('C35').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D35').Value = "'0.9999"
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').Value = "'1.111"
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.01993"
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.05313"
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = "'7.326"
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.5355"
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = "'2.873"
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.1721"
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'2.320"
$ws.Range('E43').Value = '  +15.95%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'8.681"
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = "'0.5119"
$ws.Range('E45').Value = '  -2.19%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = "'10.67"
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = "'1.698"
$ws.Range('E47').Value = '  +1.13%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'105.46"
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'1.000"
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.06410"
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = "'66.02"
$ws.Range('E51').Value = '  +4.01%  '
